# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.394.10"
$ws.Range("E2").Value = "  +2.63%  "
$ws.Range("D3").Value = "2.696.73"
$ws.Range("E3").Value = "  +1.43%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.35"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.33"
$ws.Range("E6").Value = "  +2.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.578"
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("D9").Value = "2.714.95"
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.51"
$ws.Range("E10").Value = "  +4.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.107"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("D14").Value = "3.170.57"
$ws.Range("E14").Value = "  +1.16%  "
$ws.Range("D15").Value = "60.399.11"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").Value = "2.847.53"
$ws.Range("E16").Value = "  +5.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "21.47"
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "352.58"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.56"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.63"
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.37"
$ws.Range("E22").Value = "  +3.69%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.16"
$ws.Range("E24").Value = "  +3.63%  "
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("E26").Value = "  +5.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.996"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0824"
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.40"
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("E30").Value = "  +6.94%  "
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("E32").Value = "  +1.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.20"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "147.74"
$ws.Range("E34").Value = "  -1.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.33"
$ws.Range("E35").Value = "  +7.03%  "
$ws.Range("E36").Value = "  +9.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.954"
$ws.Range("E37").Value = "  -5.65%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.880"
$ws.Range("E38").Value = "  +4.29%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.53"
$ws.Range("E39").Value = "  +9.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.96"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.72"
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "287.52"
$ws.Range("E42").Value = "  +3.28%  "
$ws.Range("E43").Value = "  +1.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.615"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0993"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("D46").Value = "2.148.46"
$ws.Range("E46").Value = "  +7.06%  "
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.92"
$ws.Range("E48").Value = "  +4.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0540"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0236"
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.40"
$ws.Range("E51").Value = "  +6.37%  "
